# Regenerate the "K" column (column G) values for the save_data sheet.
# The commit replaces the old Strike#-derived values in column G with the
# new K values (std/mean based s_vals calc), rows 2-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(7,4,6,5,4,7,1,10,7,7,7,4,2,8,10,15,5,6,5,7,3,6,12,7,3,7,7,7,10,6,6,4,12,6,6,2,2,1)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
